$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $savedStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $savedStyle
}

Set-TextValue 'D2' '26.845.17'
Set-TextValue 'E2' '  -1.17%  '

Set-TextValue 'D3' '1.559.81'
Set-TextValue 'E3' '  -0.58%  '

Set-TextValue 'E4' '  +0.04%  '

Set-TextValue 'D5' '205.35'
Set-TextValue 'E5' '  -0.77%  '

Set-TextValue 'E6' '  -1.68%  '

Set-TextValue 'E7' '  +0.06%  '

Set-TextValue 'B8' 'Solana'
Set-TextValue 'C8' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 'D8' '21.57'
Set-TextValue 'E8' '  -3.53%  '

Set-TextValue 'B9' 'Cardano'
Set-TextValue 'C9' 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue 'D9' '0.247'
Set-TextValue 'E9' '  -0.23%  '

Set-TextValue 'E10' '  -1.00%  '

Set-TextValue 'D11' '0.0861'
Set-TextValue 'E11' '  -0.45%  '

Set-TextValue 'D12' '1.781.76'
Set-TextValue 'E12' '  -0.47%  '

Set-TextValue 'D13' '1.555.45'
Set-TextValue 'E13' '  -0.42%  '

Set-TextValue 'E14' '  -1.37%  '

Set-TextValue 'E15' '  -1.31%  '

Set-TextValue 'D16' '26.841.65'
Set-TextValue 'E16' '  -1.18%  '

Set-TextValue 'D17' '61.32'
Set-TextValue 'E17' '  -2.62%  '

Set-TextValue 'D18' '214.38'
Set-TextValue 'E18' '  -0.24%  '

Set-TextValue 'D19' '7.31'
Set-TextValue 'E19' '  +0.22%  '

Set-TextValue 'D20' '0.0₃0683'
Set-TextValue 'E20' '  -0.47%  '

Set-TextValue 'E21' '  +0.03%  '

Set-TextValue 'E22' '  -0.31%  '

Set-TextValue 'D23' '9.13'
Set-TextValue 'E23' '  -2.74%  '

Set-TextValue 'E24' '  +0.69%  '

Set-TextValue 'D25' '153.16'
Set-TextValue 'E25' '  +0.57%  '

Set-TextValue 'D26' '6.58'
Set-TextValue 'E26' '  -1.21%  '

Set-TextValue 'D27' '14.94'
Set-TextValue 'E27' '  -0.11%  '

Set-TextValue 'E28' '  +0.07%  '

Set-TextValue 'E29' '  -1.55%  '

Set-TextValue 'E30' '  +0.47%  '

Set-TextValue 'E31' '  -2.71%  '

Set-TextValue 'E32' '  +0.51%  '

Set-TextValue 'D33' '1.375.59'
Set-TextValue 'E33' '  -1.53%  '

Set-TextValue 'E34' '  -0.09%  '

Set-TextValue 'E35' '  -2.87%  '

Set-TextValue 'E36' '  -0.58%  '

Set-TextValue 'D37' '0.923'
Set-TextValue 'E37' '  -2.07%  '

Set-TextValue 'E38' '  -1.55%  '

Set-TextValue 'E39' '  +1.46%  '

Set-TextValue 'D40' '0.808'
Set-TextValue 'E40' '  -0.93%  '

Set-TextValue 'E41' '  +0.06%  '

Set-TextValue 'B42' 'FraxShare'
Set-TextValue 'C42' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D42' '5.56'
Set-TextValue 'E42' '  +4.14%  '

Set-TextValue 'B43' 'WEMIXToken'
Set-TextValue 'C43' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D43' '0.990'
Set-TextValue 'E43' '  -0.21%  '

Set-TextValue 'B44' 'MXToken'
Set-TextValue 'C44' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D44' '2.18'
Set-TextValue 'E44' '  +0.52%  '

Set-TextValue 'B45' 'RenderToken'
Set-TextValue 'C45' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D45' '1.77'
Set-TextValue 'E45' '  -1.56%  '

Set-TextValue 'D46' '63.48'
Set-TextValue 'E46' '  -0.23%  '

Set-TextValue 'B47' 'mCoin'
Set-TextValue 'C47' 'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin'
Set-TextValue 'D47' '2.31'
Set-TextValue 'E47' '  -3.02%  '

Set-TextValue 'B48' 'RocketPoolETH'
Set-TextValue 'C48' 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue 'D48' '1.695.18'
Set-TextValue 'E48' '  -0.33%  '

Set-TextValue 'B49' 'Quant'
Set-TextValue 'C49' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D49' '86.47'
Set-TextValue 'E49' '  +0.74%  '

Set-TextValue 'B50' 'Cronos'
Set-TextValue 'C50' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D50' '0.0510'
Set-TextValue 'E50' '  +3.44%  '

Set-TextValue 'B51' 'BabyDogeCoin'
Set-TextValue 'C51' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D51' '0.0₇0978'
Set-TextValue 'E51' '  -0.12%  '
